# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6637336
$ws.Range("C4").Value = 1089
$ws.Range("D4").Value = 3918492
$ws.Range("E4").Value = 2521409
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 197435

# --- Row 15: Iran ---
$ws.Range("B15").Value = 399940
$ws.Range("C15").Value = 2139
$ws.Range("D15").Value = 344516
$ws.Range("E15").Value = 32395
$ws.Range("G15").Value = 116
$ws.Range("H15").Value = 23029

# --- Row 33: Kazajistan ---
$ws.Range("D33").Value = 100518
$ws.Range("E33").Value = 4577

# --- Rows 34-35: Rumania overtakes Republica Dominicana ---
$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 102386
$ws.Range("C34").Value = 1311
$ws.Range("D34").Value = 42811
$ws.Range("E34").Value = 55448
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 27
$ws.Range("H34").Value = 4127

$ws.Range("A35").Value = "Republica Dominicana"
$ws.Range("B35").Value = 102232
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 75747
$ws.Range("E35").Value = 24544
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1941

# --- Row 46: Emiratos Arabes Unidos ---
$ws.Range("B46").Value = 78849
$ws.Range("C46").Value = 1007
$ws.Range("D46").Value = 68983
$ws.Range("E46").Value = 9467
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 399

# --- Rows 60-61: Suiza overtakes Uzbekistan ---
$ws.Range("A60").Value = "Suiza"
$ws.Range("B60").Value = 46704
$ws.Range("C60").Value = 465
$ws.Range("D60").Value = 38500
$ws.Range("E60").Value = 6184
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 2020

$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 46498
$ws.Range("C61").Value = 338
$ws.Range("D61").Value = 43023
$ws.Range("E61").Value = 3093
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 382

# --- Row 85: Madagascar ---
$ws.Range("B85").Value = 15737
$ws.Range("C85").Value = 68
$ws.Range("D85").Value = 14349
$ws.Range("E85").Value = 1178
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 210

# --- Rows 146-147: Malta overtakes Botsuana ---
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2274
$ws.Range("C146").Value = 27
$ws.Range("D146").Value = 1850
$ws.Range("E146").Value = 409
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 15

$ws.Range("A147").Value = "Botsuana"
$ws.Range("B147").Value = 2252
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 546
$ws.Range("E147").Value = 1696
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 10

# --- Row 183: Gibraltar ---
$ws.Range("B183").Value = 327
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 294
$ws.Range("E183").Value = 33

# --- Update timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 13:03"
